$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# Title (appears twice: Heading1 and the bold run near the end)
Replace-Text "Play Diamond Duke for Free: Quickspin's Classic Bar Slot" "Play Diamond Duke Free - Classic Bar Slot Game"

# "What we like" bullet points
Replace-Text "Eye-catching graphics with a modern twist" "Refined and eye-catching graphics"
Replace-Text "Special symbols such as Wilds and Multipliers" "Special symbols for increased winning chances"
Replace-Text "Developed by award-winning provider Quickspin" "Similar games available for players to try"

# "What we don't like" bullet points
Replace-Text "Limited bonus features compared to other modern slots" "Limited game grid (3x3)"
Replace-Text "Smaller game grid than most video slots" "Lacks unique bonus features"

# Meta description (italic run)
Replace-Text "Read our review of Diamond Duke, the classic bar-style online slot by Quickspin. Play for free and win big with Wilds, Multipliers, and more." "Review of Diamond Duke, a classic bar slot game with refined graphics. Play for free!"
